# "Prueba Cambio en el archivo"
# The only real content edit is typing "." into cell B31 of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B31").Value = "."

# Leave the selection on the header row, matching the saved view state.
$ws.Range("A1:B1").Select()
